$wb = $excel.ActiveWorkbook

# Update the Date value on the Metadata sheet
$metaSheet = $wb.Worksheets.Item("Metadata")
$metaSheet.Range("B8").Value = "2025-05-20T15:08:54+00:00"

# Update the concept codes/displays on the Concepts sheet
$conceptsSheet = $wb.Worksheets.Item("Concepts")
$conceptsSheet.Range("B2").Value = "Yes"
$conceptsSheet.Range("C2").Value = "Yes"
$conceptsSheet.Range("B3").Value = "No"
$conceptsSheet.Range("C3").Value = "No"
$conceptsSheet.Range("B4").Value = "Unknown"
$conceptsSheet.Range("C4").Value = "Unknown"
